# Applies the two changes described by the diff:
#  1. The cached "datetimeFigureOut" field text ("2024/5/21" -> "2024/8/15")
#     on the slide master and every slide layout's date placeholder.
#  2. In the chi-square statistics textbox on slide 2, the capital "P"
#     (the italic p-value variable) is retyped as lowercase "p".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date field text wherever it reads "2024/5/21".
#    This lives on the slide master plus every custom (slide) layout.
# ---------------------------------------------------------------------
$oldDate = "2024/5/21"
$newDate = "2024/8/15"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Fix the p-value label in the SEM fit-statistics textbox on slide 2
#    ("...=25.899, P=0.055, df=16...") -> lowercase "p".
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(2)
$statsShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "*P=0.055*") {
            $statsShape = $sh
        }
    }
}

if ($statsShape -ne $null) {
    $tr = $statsShape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf("P=0.055")
    if ($idx -ge 0) {
        $charRange = $tr.Characters($idx + 1, 1)
        $charRange.Text = "p"
    }
}
